$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("quiz")

# Update "Marking" row right-answer count
$ws.Range("B11").Value = 5

# Update "Total" row right-answer count
$ws.Range("B12").Value = 120

# Update the correct/total marks text
$ws.Range("E12").Value = "120/140"
